# Add a new worksheet "数学" (Math) after "动态规划" (Dynamic Programming),
# mirroring the existing problem-list sheets, with a single "roman to int" style
# row-2 entry (time/space complexity O(n) filled in; other columns left blank).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new sheet right after the 3rd sheet (动态规划) so it lands in
#    position 4, matching sheetId=4 / r:id=rId4 in the target workbook.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Add([Type]::Missing, $ws3)
$ws4.Name = "数学"

# ---------------------------------------------------------------------------
# 2. Column widths (character units, adjusted for this host's 7px Maximum
#    Digit Width so the saved OOXML "width" lands on the intended value).
# ---------------------------------------------------------------------------
$ws4.Columns.Item(1).ColumnWidth = 6.714285714285714   # -> 7.5
$ws4.Columns.Item(2).ColumnWidth = 12.285714285714286  # -> 13
$ws4.Columns.Item(3).ColumnWidth = 64.71428571428571   # -> 65.5
$ws4.Columns.Item(4).ColumnWidth = 58.142857142857146  # -> 58.83203125
$ws4.Columns.Item(5).ColumnWidth = 25.857142857142858  # -> 26.5
$ws4.Columns.Item(6).ColumnWidth = 19.857142857142858  # -> 20.5
$ws4.Columns.Item(7).ColumnWidth = 20.142857142857142  # -> 20.83203125

# ---------------------------------------------------------------------------
# 3. Header row (row 1) - same look as the other sheets: bold 微软雅黑 14 for
#    A1:G1, plain Calibri 14 for the trailing H1:K1 filler cells.
# ---------------------------------------------------------------------------
$headerRange = $ws4.Range("A1:G1")
$headerRange.Font.Name = "微软雅黑"
$headerRange.Font.Size = 14
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4131
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true

$ws4.Range("A1").Value = "No."
$ws4.Range("B1").Value = "leetcode"
$ws4.Range("C1").Value = "题目"
$ws4.Range("D1").Value = "解题方法"
$ws4.Range("E1").Value = "解题关键词"
$ws4.Range("F1").Value = "时间复杂度"
$ws4.Range("G1").Value = "空间复杂度"

$fillerRange = $ws4.Range("H1:K1")
$fillerRange.Font.Name = "Calibri"
$fillerRange.Font.Size = 14
$fillerRange.Font.Bold = $false
$fillerRange.HorizontalAlignment = -4131
$fillerRange.VerticalAlignment = -4108
$fillerRange.WrapText = $true

$ws4.Rows.Item(1).RowHeight = 35

# ---------------------------------------------------------------------------
# 4. Data row (row 2) - matches the other sheets' data-row formatting
#    (非加粗 微软雅黑 14). Only the first / last two columns are filled in.
# ---------------------------------------------------------------------------
$dataRange = $ws4.Range("A2:G2")
$dataRange.Font.Name = "微软雅黑"
$dataRange.Font.Size = 14
$dataRange.Font.Bold = $false
$dataRange.HorizontalAlignment = -4131
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true

$ws4.Range("A2").Value = 1
$ws4.Range("F2").Value = "O(n)"
$ws4.Range("G2").Value = "O(n)"

$ws4.Rows.Item(2).RowHeight = 22

# ---------------------------------------------------------------------------
# 5. Selection / activation for the new sheet, then make it the active tab.
# ---------------------------------------------------------------------------
$ws4.Range("D14").Select()
$ws4.Activate()

# ---------------------------------------------------------------------------
# 6. 动态规划 (sheet 3) picks up new auto-computed wrap-row heights and a
#    "select rows 1:2" selection as a side effect of the edit session.
# ---------------------------------------------------------------------------
$ws3.Rows.Item(2).RowHeight = 132
$ws3.Rows.Item(3).RowHeight = 198
$ws3.Rows.Item(4).RowHeight = 308
$ws3.Range("A1:XFD2").Select()

# Re-activate the new Math sheet so it ends up as the saved active tab.
$ws4.Activate()
